# Append new task rows (166-171) to the "SB Squares Tasks" report.
# These rows document the race-condition fix for admin draft pick-on-behalf
# plus the two new future-plan docs, per the commit message / diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Task #, B=Type, C=Subject, D=Assigned To, E=Status, F=Files Changed
$rows = @(
    @{ Row = 166; A = 166; B = "feature"; C = "Gate simulation runner behind per-game simulation_enabled flag with superadmin toggle"; D = "Srini"; E = "Done"; F = "" },
    @{ Row = 167; A = 166; B = "docs";    C = "Create future plan document for server-side live score polling (pg_cron + pg_net)"; D = "Claude"; E = "Done"; F = "" },
    @{ Row = 168; A = 167; B = "docs";    C = "Create future plan document for Cloudflare Turnstile CAPTCHA on game creation"; D = "Claude"; E = "Done"; F = "docs/plan-turnstile-captcha.md" },
    @{ Row = 169; A = 168; B = "bugfix";  C = "Fix draft pick race condition: serialize pickOnBehalf with promise queue, fetch fresh squares from DB, verify pick success counts, retry on contention, add queue/active UI indicators"; D = "ui-dev"; E = "Done"; F = "" },
    @{ Row = 170; A = 169; B = "feature"; C = "Add UX feedback for draft pick operations: picking banner with spinner above draft circles, toast notifications on success/failure, disable ALL pick buttons during any active operation"; D = "ui-dev"; E = "Done"; F = "" },
    @{ Row = 171; A = 168; B = "bugfix";  C = "Fix race condition in admin draft pick-on-behalf: promise queue + fresh DB reads + verify/retry + UI queue indicators"; D = "Team (architect, uiux-expert, ui-dev)"; E = "Done"; F = "app/game/[gameId]/admin/page.tsx" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    if ($r.F -ne "") {
        $ws.Range("F$rowNum").Value = $r.F
    }
}
